$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '41.230.27'
$ws.Cells.Item(2, 5).Value = '  -1.82%  '
$ws.Cells.Item(3, 4).Value = '2.165.05'
$ws.Cells.Item(3, 5).Value = '  -2.24%  '
$ws.Cells.Item(4, 5).Value = '  +0.19%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '235.96'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -1.90%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.605'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -3.06%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '69.38'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -5.33%  '
$ws.Cells.Item(8, 5).Value = '  +0.09%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.566'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -6.65%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '38.71'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -9.20%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0907'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -4.72%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '54.32'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -5.44%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.1000'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -2.31%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.66'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -5.82%  '
$ws.Cells.Item(15, 4).Value = '2.485.40'
$ws.Cells.Item(15, 5).Value = '  -2.24%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '14.36'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +0.68%  '
$ws.Cells.Item(17, 4).Value = '2.199.50'
$ws.Cells.Item(17, 5).Value = '  -0.19%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.782'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -6.52%  '
$ws.Cells.Item(19, 4).Value = '41.030.40'
$ws.Cells.Item(19, 5).Value = '  -1.92%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0991'
$ws.Cells.Item(20, 5).Value = '  -6.92%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '69.73'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -4.45%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.77'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -6.15%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '224.47'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -2.14%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '9.29'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -14.21%  '
$ws.Cells.Item(25, 5).Value = '  +0.07%  '
$ws.Cells.Item(26, 5).Value = '  -10.25%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.62'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -9.56%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '3.37'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -6.72%  '
$ws.Cells.Item(29, 5).Value = '  -3.30%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '2.17'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -1.31%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '168.46'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +0.71%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '19.76'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -3.36%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '29.96'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +1.23%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0752'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -5.11%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.07'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -8.79%  '
$ws.Cells.Item(36, 5).Value = '  -4.05%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.102'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -7.96%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '4.07'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -4.39%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.0278'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -6.95%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.05'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -3.20%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '11.38'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -17.38%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '5.29'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -5.64%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '57.98'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -11.62%  '
$ws.Cells.Item(44, 5).Value = '  -5.54%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '8.21'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -5.25%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0958'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '96.49'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -7.38%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.07'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -4.14%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.11'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -4.83%  '
$ws.Cells.Item(50, 5).Value = '  -9.47%  '
$ws.Cells.Item(51, 5).Value = '  -3.31%  '
